# The "Glauciane" account (005981575) row is moved from its old spot
# (right after "Andre" / 005003629, row 11) to a new spot right before
# "Ahmad" / 004368468 (row 4). Its Saldo value also changes from
# 590.69 to 29390.69.
#
# Implemented as: insert a fresh blank row at row 4, fill it with the
# Glauciane data (new value), then delete the now-duplicate old
# Glauciane row (which has shifted down to row 12 because of the
# insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the "Ahmad" row (row 4).
$ws.Rows.Item(4).Insert()

# Populate the new row with the Glauciane record at its new balance.
# The account number is a zero-padded numeric string, so it must be
# written with a leading apostrophe to keep it stored as text (and not
# silently coerced to the number 5981575, losing the leading zeros).
$ws.Cells.Item(4, 1).Value = "'005981575"
$ws.Cells.Item(4, 2).Value = "Glauciane"
$ws.Cells.Item(4, 3).Value = 29390.69

# Remove the original Glauciane row, which the insert above shifted
# down from row 11 to row 12.
$ws.Rows.Item(12).Delete()
